$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '27.574.85'
    'E2' = '  -1.29%  '
    'D3' = '1.846.88'
    'E3' = '  -2.18%  '
    'D4' = '1.005'
    'E4' = '  -1.11%  '
    'D5' = '333.46'
    'E5' = '  -0.50%  '
    'D6' = '1.003'
    'E6' = '  -1.00%  '
    'D7' = '0.4642'
    'E7' = '  -1.29%  '
    'D8' = '0.3858'
    'E8' = '  -1.85%  '
    'E9' = '  -0.96%  '
    'D10' = '0.07912'
    'E10' = '  -0.85%  '
    'D11' = '0.9947'
    'E11' = '  -1.92%  '
    'D12' = '21.48'
    'E12' = '  -1.48%  '
    'D13' = '1.839.44'
    'E13' = '  -2.40%  '
    'D14' = '5.922'
    'E14' = '  -0.76%  '
    'D15' = '7.112'
    'E15' = '  -0.57%  '
    'D16' = '1.005'
    'E16' = '  -1.21%  '
    'D17' = '88.92'
    'E17' = '  +1.49%  '
    'D18' = '0.06643'
    'E18' = '  -1.57%  '
    'D19' = '0.00001035'
    'E19' = '  -1.48%  '
    'D20' = '17.07'
    'E20' = '  -0.54%  '
    'E21' = '  -1.00%  '
    'D22' = '27.572.10'
    'E22' = '  -1.34%  '
    'D23' = '5.379'
    'E23' = '  -2.19%  '
    'D24' = '10.91'
    'E24' = '  -0.59%  '
    'D25' = '2.298'
    'E25' = '  -2.95%  '
    'D26' = '158.17'
    'E26' = '  -0.52%  '
    'D27' = '19.52'
    'E27' = '  -2.43%  '
    'D28' = '2.102'
    'E28' = '  -0.01%  '
    'D29' = '5.407'
    'E29' = '  -1.80%  '
    'D30' = '119.86'
    'D31' = '0.9770'
    'E31' = '  +1.40%  '
    'D32' = '0.09408'
    'E32' = '  -1.63%  '
    'D33' = '3.584'
    'E33' = '  -1.70%  '
    'E34' = '  -1.26%  '
    'D35' = '1.341'
    'E35' = '  -1.42%  '
    'D36' = '0.06021'
    'E36' = '  -1.74%  '
    'D37' = '0.02226'
    'E37' = '  -0.98%  '
    'D38' = '8.304'
    'E38' = '  +1.22%  '
    'D39' = '1.184'
    'E39' = '  -2.46%  '
    'D40' = '0.5891'
    'D41' = '0.1865'
    'E41' = '  -1.89%  '
    'D42' = '10.30'
    'E42' = '  -0.26%  '
    'D43' = '1.242'
    'E43' = '  -1.90%  '
    'D44' = '0.5577'
    'E44' = '  -1.89%  '
    'D45' = '12.22'
    'E45' = '  -0.48%  '
    'D46' = '1.902'
    'E46' = '  -1.96%  '
    'E47' = '  -2.55%  '
    'D48' = '110.72'
    'E48' = '  -2.73%  '
    'D49' = '1.053'
    'E49' = '  -1.50%  '
    'D50' = '1.002'
    'E50' = '  -1.10%  '
    'D51' = '70.03'
    'E51' = '  -1.44%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

Write-Host "Updated $($updates.Count) cells"